# Rename the single worksheet from "Property1" to "DataNode"
# (unifying the DataNode / DataTable / Entity naming, per the commit message)
# and update the current cell selection to C38, matching the author's
# last recorded cursor position when the change was committed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "DataNode"
$ws.Range("C38").Select()
